$d = $word.ActiveDocument

# Position at the very end of the document body (right after "${fin}"),
# so the new content is appended as new paragraphs rather than replacing
# the existing last paragraph.
$insertionRange = $d.Content
$insertionRange.Collapse(0)

# Raw OOXML for the two new paragraphs appended after "${fin}":
#   1) a paragraph holding only a page break
#   2) a right-aligned paragraph with the verso marker "${textVerso}"
$newParagraphsXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:bidi w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Shlomo Stam" w:hAnsi="Shlomo Stam" w:cs="Shlomo Stam"/><w:spacing w:val="4"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Shlomo Stam" w:hAnsi="Shlomo Stam" w:cs="Shlomo Stam"/><w:spacing w:val="4"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:lang w:val="fr-FR"/></w:rPr><w:br w:type="page"/></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:suppressAutoHyphens/><w:spacing w:after="0" w:line="460" w:lineRule="exact"/><w:jc w:val="right"/><w:rPr><w:rFonts w:ascii="Shlomo Stam" w:hAnsi="Shlomo Stam" w:cs="Shlomo Stam"/><w:spacing w:val="4"/><w:sz w:val="30"/><w:szCs w:val="30"/><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Shlomo Stam" w:hAnsi="Shlomo Stam" w:cs="Shlomo Stam"/><w:spacing w:val="4"/><w:sz w:val="30"/><w:szCs w:val="30"/><w:lang w:val="fr-FR"/></w:rPr><w:lastRenderedPageBreak/><w:t>${</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Shlomo Stam" w:hAnsi="Shlomo Stam" w:cs="Shlomo Stam"/><w:spacing w:val="4"/><w:sz w:val="30"/><w:szCs w:val="30"/><w:lang w:val="fr-FR"/></w:rPr><w:t>t</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Shlomo Stam" w:hAnsi="Shlomo Stam" w:cs="Shlomo Stam"/><w:spacing w:val="4"/><w:sz w:val="30"/><w:szCs w:val="30"/><w:lang w:val="fr-FR"/></w:rPr><w:t>exteVerso</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Shlomo Stam" w:hAnsi="Shlomo Stam" w:cs="Shlomo Stam"/><w:spacing w:val="4"/><w:sz w:val="30"/><w:szCs w:val="30"/><w:lang w:val="fr-FR"/></w:rPr><w:t>}</w:t></w:r></w:p>
'@

$null = $insertionRange.InsertXML($newParagraphsXml)

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
